$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.147.15"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "1.642.17"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("D5").Value = "'216.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "'19.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "1.869.67"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "1.640.88"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "'63.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "26.152.53"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "'0.999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "'194.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").Value = "'10.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").Value = "'1.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "'142.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "'0.0502"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("D32").Value = "'3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "1.131.25"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("D39").Value = "'2.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'100.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").Value = "'0.797"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "1.778.96"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "'56.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("E48").Value = "  +4.27%  "
$ws.Range("D49").Value = "'0.0516"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +2.89%  "
